$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "Related Works" paragraph: tighten the opening sentence.
#    "Basically, there are 4 main types of fake news detection techniques viz. "
#    -> "There are 4 main types of fake news detection techniques. "
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Text = "Basically, there are 4 main types of fake news detection techniques viz. "
$rng.Find.Forward = $true
$rng.Find.Wrap = 1
$rng.Find.MatchCase = $false
if ($rng.Find.Execute()) {
    $rng.Text = "There are 4 main types of fake news detection techniques. "
}

# ---------------------------------------------------------------------------
# 2. "Experiment" paragraph: wording is unchanged, but re-create the _GoBack
#    bookmark Word drops at the point of the most recent edit before saving
#    ("...gave a certain score which|_GoBack| could then...").
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Text = "score which could"
$rng.Find.Forward = $true
$rng.Find.Wrap = 1
if ($rng.Find.Execute()) {
    $markStart = $rng.Start + "score which".Length
    $goBackRange = $d.Range($markStart, $markStart)
    $d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null
}

# ---------------------------------------------------------------------------
# 3. References: add a new bulleted hyperlink entry after the
#    towardsdatascience.com link, then a trailing (non-bulleted) blank
#    paragraph to close out the list.
# ---------------------------------------------------------------------------
$lastRef = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*towardsdatascience*") {
        $lastRef = $p
    }
}

$insertPoint = $lastRef.Range
$insertPoint.Collapse(0)
$insertPoint.InsertParagraphAfter()

$newPara = $d.Paragraphs($d.Paragraphs.Count)
$newRange = $newPara.Range
$url = "https://www.researchgate.net/profile/Sandeep-Pande-2/publication/362235353_Fake_News_Identification_Using_Regression_Analysis_and_Web_Scraping/links/62e0b4819d410c5ff367263c/Fake-News-Identification-Using-Regression-Analysis-and-Web-Scraping.pdf"
$h = $d.Hyperlinks.Add($newRange, $url, "", "", $url)
$h.Range.Font.Name = "Times New Roman"

# Trailing empty paragraph (list formatting removed, small left indent) that
# closes out the References list.
$p2 = $d.Paragraphs($d.Paragraphs.Count)
$tailRange = $p2.Range
$tailRange.Collapse(0)
$tailRange.InsertParagraphAfter()

$blankPara = $d.Paragraphs($d.Paragraphs.Count)
$blankPara.Range.ListFormat.RemoveNumbers()
$blankPara.Style = "Normal"
$blankPara.Format.LineSpacingRule = 2
$blankPara.Format.LineSpacing = 24
$blankPara.LeftIndent = 18
$blankPara.Range.Font.Name = "Times New Roman"

Write-Output "done"
